$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.07410728931427
$ws.Range("B1").Value = 4.477452754974365
$ws.Range("C1").Value = 2.459789276123047
$ws.Range("D1").Value = 1.699832916259766
$ws.Range("E1").Value = 1.378459334373474
